$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.016392150120196
$ws.Range("D2").Value = 1.022474799646343
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.014736051652402
$ws.Range("I2").Value = 1.026577868333675
$ws.Range("J2").Value = 1.021612221934027
$ws.Range("K2").Value = 1.025309072660218
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.017593311049693
$ws.Range("N2").Value = 1.011223809385837
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.017726981172331
$ws.Range("D3").Value = 1.023458789068739
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.016720047224367
$ws.Range("I3").Value = 1.026811579663227
$ws.Range("J3").Value = 1.02258045626705
$ws.Range("K3").Value = 1.026099251644517
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.019379047084176
$ws.Range("N3").Value = 1.011551500582931
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018588652532026
$ws.Range("D4").Value = 1.024093584111876
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.018001286891546
$ws.Range("I4").Value = 1.026960671422299
$ws.Range("J4").Value = 1.023204530642104
$ws.Range("K4").Value = 1.026608036569209
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.02053163716861
$ws.Range("N4").Value = 1.011762519956938
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018950415550297
$ws.Range("D5").Value = 1.024359999106261
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.018539329908409
$ws.Range("I5").Value = 1.027022840140126
$ws.Range("J5").Value = 1.023466313982192
$ws.Range("K5").Value = 1.026821332647641
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.021015507418645
$ws.Range("N5").Value = 1.011850990580579
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019011129024574
$ws.Range("D6").Value = 1.024404704941814
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.018629635597962
$ws.Range("I6").Value = 1.027033248715232
$ws.Range("J6").Value = 1.023510234867411
$ws.Range("K6").Value = 1.026857111118439
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.021096712023836
$ws.Range("N6").Value = 1.011865831057976
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018593488317643
$ws.Range("D7").Value = 1.024097145735476
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.018008478546043
$ws.Range("I7").Value = 1.026961504123518
$ws.Range("J7").Value = 1.023208030865691
$ws.Range("K7").Value = 1.026610888982343
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.020538105313749
$ws.Range("N7").Value = 1.011763703053918
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.016843692576305
$ws.Range("D8").Value = 1.022807742302091
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.015407088355146
$ws.Range("I8").Value = 1.026657295281395
$ws.Range("J8").Value = 1.021939949768448
$ws.Range("K8").Value = 1.025576640664084
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.018197419011811
$ws.Range("N8").Value = 1.011334766288359
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013744239440717
$ws.Range("D9").Value = 1.020520785121757
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.010802879310382
$ws.Range("I9").Value = 1.026104811441196
$ws.Range("J9").Value = 1.019686479882374
$ws.Range("K9").Value = 1.023734698939421
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.014049905364543
$ws.Range("N9").Value = 1.01057103351803
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.011666610473976
$ws.Range("D10").Value = 1.018985840970813
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.007718700992413
$ws.Range("I10").Value = 1.025725333837326
$ws.Range("J10").Value = 1.01817103555745
$ws.Range("K10").Value = 1.022493351258259
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.011268484857558
$ws.Range("N10").Value = 1.01005644664783
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.010764173176352
$ws.Range("D11").Value = 1.018318677980603
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.006379482515291
$ws.Range("I11").Value = 1.025558345156614
$ws.Range("J11").Value = 1.017511628882506
$ws.Range("K11").Value = 1.021952591070258
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.010059983366351
$ws.Range("N11").Value = 1.009832308057868
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01042853560276
$ws.Range("D12").Value = 1.01807047881044
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.005881452527354
$ws.Range("I12").Value = 1.025495914415651
$ws.Range("J12").Value = 1.017266206027127
$ws.Range("K12").Value = 1.021751234655097
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.009610452592411
$ws.Range("N12").Value = 1.009748852328655
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.010500550720009
$ws.Range("D13").Value = 1.018123735903677
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.005988308451405
$ws.Range("I13").Value = 1.025509324326133
$ws.Range("J13").Value = 1.017318872354307
$ws.Range("K13").Value = 1.021794448740903
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.009706907721089
$ws.Range("N13").Value = 1.009766762987414
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.010736438144252
$ws.Range("D14").Value = 1.018298169652596
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.006338327251715
$ws.Range("I14").Value = 1.025553192860577
$ws.Range("J14").Value = 1.017491352190502
$ws.Range("K14").Value = 1.021935956998458
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.010022838170727
$ws.Range("N14").Value = 1.009825413689589
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.010881718603259
$ws.Range("D15").Value = 1.018405592821326
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.006553907275334
$ws.Range("I15").Value = 1.02558016815021
$ws.Range("J15").Value = 1.01759755756755
$ws.Range("K15").Value = 1.022023079266329
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.010217408024597
$ws.Range("N15").Value = 1.009861523655898
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011726442208004
$ws.Range("D16").Value = 1.019030064715479
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.00780749970702
$ws.Range("I16").Value = 1.025736359833827
$ws.Range("J16").Value = 1.018214729935931
$ws.Range("K16").Value = 1.022529170782805
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.011348600504121
$ws.Range("N16").Value = 1.010071293983687
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012255555619857
$ws.Range("D17").Value = 1.019421099846687
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.008592826407948
$ws.Range("I17").Value = 1.025833617725324
$ws.Range("J17").Value = 1.018601001235489
$ws.Range("K17").Value = 1.022845754733325
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.012057048893431
$ws.Range("N17").Value = 1.010202522457996
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012563908053096
$ws.Range("D18").Value = 1.019648941107469
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.009050533926574
$ws.Range("I18").Value = 1.025890088914216
$ws.Range("J18").Value = 1.018825997625376
$ws.Range("K18").Value = 1.023030099774288
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.01246987787098
$ws.Range("N18").Value = 1.010278938705583
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012669002688034
$ws.Range("D19").Value = 1.019726588073484
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.009206539899408
$ws.Range("I19").Value = 1.025909300487392
$ws.Range("J19").Value = 1.018902663378409
$ws.Range("K19").Value = 1.023092903754491
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.012610575130036
$ws.Range("N19").Value = 1.01030497318773
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.012198814806111
$ws.Range("D20").Value = 1.019379170648707
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.008508605702822
$ws.Range("I20").Value = 1.025823209544503
$ws.Range("J20").Value = 1.018559590010456
$ws.Range("K20").Value = 1.022811820691898
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.011981080287058
$ws.Range("N20").Value = 1.01018845604939
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.010666987188087
$ws.Range("D21").Value = 1.018246813943035
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.006235271750239
$ws.Range("I21").Value = 1.025540285832892
$ws.Range("J21").Value = 1.017440574762398
$ws.Range("K21").Value = 1.021894300028788
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.009929822337541
$ws.Range("N21").Value = 1.009808148079539
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.009701360884009
$ws.Range("D22").Value = 1.017532624953318
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.004802543284303
$ws.Range("I22").Value = 1.025360063605798
$ws.Range("J22").Value = 1.01673416742502
$ws.Range("K22").Value = 1.021314557208465
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.008636405305711
$ws.Range("N22").Value = 1.009567871046707
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.010213498506879
$ws.Range("D23").Value = 1.01791144365724
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.005562388588955
$ws.Range("I23").Value = 1.025455825073302
$ws.Range("J23").Value = 1.017108918927053
$ws.Range("K23").Value = 1.021622163061496
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.00932242829356
$ws.Range("N23").Value = 1.00969535746525
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.01222445437329
$ws.Range("D24").Value = 1.019398117399222
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.008546662527033
$ws.Range("I24").Value = 1.025827913349718
$ws.Range("J24").Value = 1.018578302915774
$ws.Range("K24").Value = 1.022827154991872
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.01201540845658
$ws.Range("N24").Value = 1.010194812446735
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.014547481142575
$ws.Range("D25").Value = 1.021113811709854
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.011995687890799
$ws.Range("I25").Value = 1.026249599812177
$ws.Range("J25").Value = 1.020271341369031
$ws.Range("K25").Value = 1.024213222263916
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.015124955898322
$ws.Range("N25").Value = 1.010769424784512
